# metadata.xlsx: reclassify the "Condition" column (E) so that the generic
# "NA" placeholder is replaced with the real, informative condition label.
# Rows whose Treatment (D) is CpG and whose Condition was "NA" become
# "Control"; rows whose Treatment is GpC and whose Condition was "NA" become
# "Sham". Rows already labelled "Burn" are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 80 }

for ($row = 2; $row -le $lastRow; $row++) {
    $treatment = $ws.Cells.Item($row, 4).Value()
    $condition = $ws.Cells.Item($row, 5).Value()

    if ($condition -eq "NA") {
        if ($treatment -eq "CpG") {
            $ws.Cells.Item($row, 5).Value = "Control"
        } elseif ($treatment -eq "GpC") {
            $ws.Cells.Item($row, 5).Value = "Sham"
        }
    }
}

$ws.Range("E80").Select() | Out-Null
